$d = $word.ActiveDocument

# Replacement 1
$ok1 = $d.Content.Find.Execute("             business/sophisticated", $false, $false, $false, $false, $false, $true, 1, $false, "business/sophisticated", 2)
if (-not $ok1) { Write-Output "FAILED replacement 1" }

# Replacement 2
$ok2 = $d.Content.Find.Execute("International trading laws are very ………………….Compliance with them is", $false, $false, $false, $false, $false, $true, 1, $false, "International trading laws are very …stringent……………….Compliance with them is", 2)
if (-not $ok2) { Write-Output "FAILED replacement 2" }

# Replacement 3
$ok3 = $d.Content.Find.Execute("………………………………..number of our workers are competent", $false, $false, $false, $false, $false, $true, 1, $false, "…Prevalent……………………………..number of our workers are competent", 2)
if (-not $ok3) { Write-Output "FAILED replacement 3" }

# Replacement 4
$ok4 = $d.Content.Find.Execute("Our views are completely …………………………………We don’t have any misunderstandings", $false, $false, $false, $false, $false, $true, 1, $false, "Our views are completely …coherent………………………………We don’t have any misunderstandings", 2)
if (-not $ok4) { Write-Output "FAILED replacement 4" }

# Replacement 5
$ok5 = $d.Content.Find.Execute("Current situation in Poland is not ……………………….to do any business", $false, $false, $false, $false, $false, $true, 1, $false, "Current situation in Poland is not …conducive…………………….to do any business", 2)
if (-not $ok5) { Write-Output "FAILED replacement 5" }

# Replacement 6
$ok6 = $d.Content.Find.Execute("Thoughtless decisions of management board have ……………………………..influence on our liquidity", $false, $false, $false, $false, $false, $true, 1, $false, "Thoughtless decisions of management board have ……disruptive………………………..influence on our liquidity", 2)
if (-not $ok6) { Write-Output "FAILED replacement 6" }

# Replacement 7
$ok7 = $d.Content.Find.Execute("In Poland there has been ………………………….trend to invest in funds", $false, $false, $false, $false, $false, $true, 1, $false, "In Poland there has been ……attentive…………………….trend to invest in funds", 2)
if (-not $ok7) { Write-Output "FAILED replacement 7" }

# Replacement 8
$ok8 = $d.Content.Find.Execute("My boss is extremelly …………………………to details . He is so meticulous", $false, $false, $false, $false, $false, $true, 1, $false, "My boss is extremelly ……stringent……………………to details . He is so meticulous", 2)
if (-not $ok8) { Write-Output "FAILED replacement 8" }

# Replacement 9
$ok9 = $d.Content.Find.Execute("Most of my clients are …………………………….of help and devotion we provide tchem with", $false, $false, $false, $false, $false, $true, 1, $false, "Most of my clients are ……appreciative……………………….of help and devotion we provide them with", 2)
if (-not $ok9) { Write-Output "FAILED replacement 9" }

# Replacement 10
$ok10 = $d.Content.Find.Execute("During negotiations some signs of courtesy and kindness can be ………………………..", $false, $false, $false, $false, $false, $true, 1, $false, "During negotiations some signs of courtesy and kindness can be …deceptive……………………..", 2)
if (-not $ok10) { Write-Output "FAILED replacement 10" }

# Replacement 11
$ok11 = $d.Content.Find.Execute("Our services and products are………………………… in all type of companies", $false, $false, $false, $false, $false, $true, 1, $false, "Our services and products are………………………… in alltype of companies", 2)
if (-not $ok11) { Write-Output "FAILED replacement 11" }

# Replacement 12
$ok12 = $d.Content.Find.Execute("Everyday I deal with the same papers , documents. My work is extremely………………………….. ", $false, $false, $false, $false, $false, $true, 1, $false, "Everyday I deal with the same papers , documents. My workisextremely………………………….. ", 2)
if (-not $ok12) { Write-Output "FAILED replacement 12" }

# Replacement 13
$ok13 = $d.Content.Find.Execute(" Inexcusable, Inevitable, Immeasurable, unattainable, Inestimable, unjustifiable, applicable, undeniable, sustainable, reconcilable", $false, $false, $false, $false, $false, $true, 1, $false, "Inexcusable, Inevitable, Immeasurable, unattainable, Inestimable, unjustifiable, applicable, undeniable, sustainable, reconcilable", 2)
if (-not $ok13) { Write-Output "FAILED replacement 13" }

# Replacement 14
$ok14 = $d.Content.Find.Execute("Some people are just too ………………………We appreciate devotion and commitment however trying to be everywhere and doing everything is not perceived as suitable behaviour ", $false, $false, $false, $false, $false, $true, 1, $false, "Some people are just too ………………………We appreciate devotion and commitment however trying to be everywhere and doing everything is not perceived as suitable behaviour", 2)
if (-not $ok14) { Write-Output "FAILED replacement 14" }

